{"js": "// Fix the typo in the cover page author name: \"Caronlina Scudeler\" -> \"Carolina Scudeler\"\n// (commit: \"Acerto no nome da Carolina na capa\")\n\nconst body = context.document.body;\n\n// Locate the misspelled name on the cover page.\nconst misspelled = \"Caronlina Scudeler\";\nconst corrected = \"Carolina Scudeler\";\n\nconst results = body.search(misspelled, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the text in place; Word keeps the run's existing formatting\n  // (Times New Roman, bold, 28/28) because insertText(\"Replace\") rewrites\n  // the matched range's text only.\n  results.items[0].insertText(corrected, \"Replace\");\n  await context.sync();\n} else {\n  // Idempotency / defensive fallback: if the exact phrase was not found\n  // (e.g. the document was already fixed, or only the misspelled token\n  // itself is present), try just the misspelled word on its own.\n  const wordResults = body.search(\"Caronlina\", { matchCase: true, matchWholeWord: false });\n  wordResults.load(\"items\");\n  await context.sync();\n  if (wordResults.items.length > 0) {\n    wordResults.items[0].insertText(\"Carolina\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Fix the typo in the cover page author name: \"Caronlina Scudeler\" -> \"Carolina Scudeler\"\n# (commit: \"Acerto no nome da Carolina na capa\")\n\n$d = $word.ActiveDocument\n\n$misspelled = \"Caronlina Scudeler\"\n$corrected  = \"Carolina Scudeler\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $misspelled\n$find.Replacement.Text = $corrected\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$found = $find.Execute($misspelled, $false, $false, $false, $false, $false, $true, 1, $false, $corrected, 2)\n\nif (-not $found) {\n    # Idempotency / defensive fallback: only the misspelled token itself.\n    $find2 = $d.Content.Find\n    $find2.Text = \"Caronlina\"\n    $find2.Replacement.Text = \"Carolina\"\n    $find2.Execute(\"Caronlina\", $false, $false, $false, $false, $false, $true, 1, $false, \"Carolina\", 2) | Out-Null\n}\n"}
